$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix mislabeled / duplicated command descriptions in the serial ICD table.
# New unique text is introduced in this order so the shared-string table
# is built up the same way the original authoring tool produced it.
$ws.Range("C17").Value = "waypoints"
$ws.Range("F17").Value = "read waypoints"

$ws.Range("C9").Value  = "transmitter slope cal"
$ws.Range("C10").Value = "transmitter offset cal"

$ws.Range("F9").Value  = "read transmitter slope values"
$ws.Range("F10").Value = "read transmitter offset values"

$ws.Range("C18").Value = "camera values"
$ws.Range("F18").Value = "read camera values"

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("F10").Select()
